$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" updates ---
$ws1 = $wb.Worksheets.Item("Schedule")
$ws1.Range("E3").Value = 583.6327769999999
$ws1.Range("F3").Value = 17.15557839506173
$ws1.Range("B4").Value = 46038.125
$ws1.Range("C4").Value = 5
$ws1.Range("D4").Value = 18.9
$ws1.Range("E4").Value = 543.31524975
$ws1.Range("F4").Value = 28.74683861111112
$ws1.Range("A5").Value = 46038.29166666666
$ws1.Range("C5").Value = 9
$ws1.Range("D5").Value = 34.02
$ws1.Range("E5").Value = 583.0842810000001
$ws1.Range("F5").Value = 17.13945564373898

# --- Sheet "Detailed" updates ---
$ws2 = $wb.Worksheets.Item("Detailed")
$ws2.Range("B31").Value = 52.15699
$ws2.Range("B32").Value = 58.73807
$ws2.Range("B33").Value = 42.97119
$ws2.Range("B34").Value = 43.22303
$ws2.Range("C34").Value = "historical"
$ws2.Range("B35").Value = 43.7905
$ws2.Range("B36").Value = 55.05363
$ws2.Range("B37").Value = 17.535
$ws2.Range("B38").Value = 57.89137
$ws2.Range("B39").Value = 62.83836
$ws2.Range("B40").Value = 80.02
$ws2.Range("B41").Value = 120.01
$ws2.Range("B42").Value = 129.81998
$ws2.Range("B43").Value = 120.01
$ws2.Range("B44").Value = 90.2877
$ws2.Range("B46").Value = 68.2715
$ws2.Range("B49").Value = 58.42568
$ws2.Range("B51").Value = 57.88079
$ws2.Range("B52").Value = 59.1159
$ws2.Range("B54").Value = 36.07
$ws2.Range("B55").Value = 49.23254
$ws2.Range("E55").Value = "ON"
$ws2.Range("B56").Value = 36.07
$ws2.Range("B57").Value = 36.07
$ws2.Range("B58").Value = 56.98
$ws2.Range("B59").Value = 61.77255
$ws2.Range("B60").Value = 62.9539
$ws2.Range("B61").Value = 73.20005
$ws2.Range("B62").Value = 60.8552
$ws2.Range("B63").Value = 57.06003
$ws2.Range("E63").Value = "OFF"
$ws2.Range("B64").Value = 36.06
$ws2.Range("B66").Value = 36.05989
$ws2.Range("B68").Value = 45.42611
$ws2.Range("B70").Value = 45.51652
$ws2.Range("B71").Value = 36.05922
$ws2.Range("B72").Value = 36.05989
$ws2.Range("B77").Value = 36.05952
$ws2.Range("B78").Value = 36.0601
$ws2.Range("B79").Value = 8.562340000000001
$ws2.Range("B80").Value = 12.68053
$ws2.Range("B81").Value = 17.07084
$ws2.Range("B82").Value = 17.87116
$ws2.Range("B83").Value = 0.34135
$ws2.Range("B84").Value = -9.3123
$ws2.Range("B85").Value = -9.55317
$ws2.Range("B86").Value = -6
$ws2.Range("B87").Value = -5.99309
$ws2.Range("B88").Value = -3.0714
$ws2.Range("B89").Value = 32.40461
$ws2.Range("B90").Value = 32.40461
$ws2.Range("B91").Value = 32.40461
$ws2.Range("B92").Value = 32.40461
$ws2.Range("B93").Value = 78
$ws2.Range("B94").Value = 64.8901
$ws2.Range("B95").Value = 57.09
$ws2.Range("B96").Value = 57.09
$ws2.Range("B97").Value = 57.06003
